# Semana 4 update: correct a handful of typos in the task labels and
# scroll/reselect the sheet as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "LOCS/HORA"
$ws.Range("A24").Value = "Crear una etiqueta de referencia para cada línea agregada o borrada para indicar el número del cambio"
$ws.Range("A33").Value = "En la etiqueta del programa indicar las líneas agregadas, borradas y totales de todo el programa"
$ws.Range("A38").Value = "Para líneas muy grandes pasar a la siguiente línea"
$ws.Range("A53").Value = "Crear fábricas para permitir el funcionamiento en diversos lenguajes"

$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("A58").Select()
